# Slide 1, Title placeholder: "Types of Health Insurance"
#   -> "Types of Commerical Health Insurance" (split into 4 runs so the
#      misspelled word "Commerical" stands alone, matching the OOXML diff)
# Also turns on "Shrink text on overflow" (a:normAutofit) for that title.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# New full title text.
$tr.Text = "Types of Commerical Health Insurance"

# Shrink text on overflow -> <a:bodyPr><a:normAutofit/></a:bodyPr>
$tf.AutoSize = 2

# Re-assigning a sub-range's own text (without changing it) forces the
# paragraph to split into separate <a:r> runs at that boundary without
# introducing any extra direct-formatting attributes.

# Run 1: "Types of "   (chars 1-9)
$r1 = $tr.Characters(1, 9)
$r1.Text = $r1.Text

# Run 2: "Commerical"  (chars 10-19)
$r2 = $tr.Characters(10, 10)
$r2.Text = $r2.Text

# Run 3: " Health "    (chars 20-27)
$r3 = $tr.Characters(20, 8)
$r3.Text = $r3.Text

# Run 4: "Insurance"   (chars 28-36)
$r4 = $tr.Characters(28, 9)
$r4.Text = $r4.Text
